$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 13-15 (the three "docente" rows that held only B/C values)
$ws.Rows("13:15").Delete()

# Update cell contents that changed independently of the row shift
$ws.Range("B10").Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range("C10").Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range("B13").Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range("C13").Value = '7459752 - Maria Ismenia Sodero Toledo Faria'
$ws.Range("B15").Value = '2166002 - Sandra Giacomin Schneider'
$ws.Range("C15").Value = '2166002 - Sandra Giacomin Schneider'
$ws.Range("B18").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("C18").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("B19").Value = 'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.'
$ws.Range("C19").Value = 'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.'
$ws.Range("B20").Value = 'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.'
$ws.Range("C20").Value = 'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.'
$ws.Range("B21").Value = 'não tem norma de recuperação'
$ws.Range("C21").Value = 'não tem norma de recuperação'
